$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Tasks" worksheet right after "Android studio " (position 2)
# ---------------------------------------------------------------------------
$android = $wb.Worksheets.Item(1)
$firebase = $wb.Worksheets.Item(2)

$tasks = $wb.Worksheets.Add($firebase)
$tasks.Name = "Tasks"

# Header row for the Tasks sheet
$tasks.Range("A1").Value = "Pending Tasks"
$tasks.Range("B1").Value = "Floating Tasks"
$tasks.Range("C1").Value = "Developer name"
$tasks.Range("D1").Value = "status"
$tasks.Range("E1").Value = "Error status "
$tasks.Range("F1").Value = "Descprition"

# Column widths (best effort - engine stores widths at 1/6-character granularity)
$tasks.Columns.Item(1).ColumnWidth = 11.619791666666666
$tasks.Columns.Item(2).ColumnWidth = 11.529947916666666
$tasks.Columns.Item(3).ColumnWidth = 13.799479166666666
$tasks.Columns.Item(5).ColumnWidth = 10.166666666666666
$tasks.Columns.Item(6).ColumnWidth = 9.436197916666666

# ---------------------------------------------------------------------------
# 2. Populate the (previously empty) "Firebase" sheet
# ---------------------------------------------------------------------------
$firebase = $wb.Worksheets.Item("Firebase")
$firebase.Range("A1").Value = "Configuration settings/errors faced"
$firebase.Range("A2").Value = "Google sign in we need to add SHA1 finger print in the project settings to enable G sign up."

$firebase.Columns.Item(1).ColumnWidth = 76.07291666666667
$firebase.Columns.Item(2).ColumnWidth = 16.529947916666668

# ---------------------------------------------------------------------------
# 3. Populate the (previously empty) "Gradle" sheet
# ---------------------------------------------------------------------------
$gradle = $wb.Worksheets.Item("Gradle")
$gradle.Range("A1").Value = "Gradle Dependencies to note"
$gradle.Range("B1").Value = "Usage purpose"

$gradle.Range("A2").Value = "`n    implementation 'com.fasterxml.jackson.core:jackson-core:2.10.2'`n    implementation 'com.fasterxml.jackson.core:jackson-annotations:2.10.2'`n    implementation 'com.fasterxml.jackson.core:jackson-databind:2.10.2'"
$gradle.Range("B2").Value = "For object mapper i.e to conver object to any type of collections"
$gradle.Range("A2").WrapText = $true
$gradle.Rows.Item(2).RowHeight = 64.5

$gradle.Range("A3").Value = "    implementation 'com.google.android.material:material:1.1.0'"
$gradle.Range("B3").Value = "Material design to support some ui functions such as Chip etc."

$gradle.Range("A4").Value = "    implementation 'com.google.android.gms:play-services-auth:16.0.1'"
$gradle.Range("B4").Value = "Google sign up dependencies."

$gradle.Columns.Item(1).ColumnWidth = 83.07291666666667
$gradle.Columns.Item(2).ColumnWidth = 53.619791666666664

# ---------------------------------------------------------------------------
# 4. Workbook / sheet selections (order matters: the last Select() call wins
#    the "active sheet" / tabSelected slot)
# ---------------------------------------------------------------------------
$android.Range("D6").Select() | Out-Null

$firebase.Range("A3").Select() | Out-Null

$gradle.Range("A7").Select() | Out-Null

$github = $wb.Worksheets.Item("GitHub")
$github.Range("H9").Select() | Out-Null

# Tasks becomes the active tab (selected last)
$tasks.Range("F6").Select() | Out-Null
